# "Add files via upload" -- refreshed inventory export for south_des_moines.xlsx
# (DSM-Bloomfield Shop). The source system re-exported the report later in the day
# with updated on-hand quantities, so the Total Cost (col F) and Quantity (col K)
# columns move together (Total Cost = Cost * Quantity), and the two grand-total rows
# (183 and 185) are refreshed to match the new column sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the "Exported On" timestamp (A2) -------------------------------
# Cell A2 holds "Exported On: <date>" as rich text ("Exported On: " bold, the date
# not bold). Only the date/time portion changed, so rewrite just those characters
# (positions 14-32) and leave the "Exported On: " prefix alone.
$exportedOnCell = $ws.Range("A2")
$exportedOnCell.Characters(14, 19).Text = "07/22/2025 07:23 AM"

# --- Body: refreshed Total Cost (F) / Quantity (K) per row -------------------------
$cellValues = [ordered]@{
    "F11" = 11.66
    "K11" = 11.0
    "F12" = 46.4
    "K12" = 10.0
    "F14" = 0.0
    "K14" = -1.0
    "F15" = 23.53
    "K15" = 13.0
    "F16" = 46.97
    "K16" = 7.0
    "F18" = 34.41
    "K18" = 3.0
    "F20" = 875.26
    "K20" = 2.0
    "F22" = 2139.96
    "K22" = 4.0
    "F32" = 128.28
    "K32" = 4.0
    "F34" = 26.64
    "K34" = 24.0
    "F35" = 6.24
    "K35" = 3.0
    "F38" = 91.58
    "K38" = 1.0
    "F39" = 18.36
    "K39" = 9.0
    "F40" = 4.83
    "K40" = 3.0
    "F45" = 438.5
    "K45" = 10.0
    "F49" = 352.65
    "K49" = 15.0
    "F50" = 69.72
    "K50" = 28.0
    "F54" = 9.6
    "K54" = 5.0
    "F55" = 30.16
    "K55" = 13.0
    "F57" = 37.35
    "K57" = 5.0
    "F59" = 36.76
    "K59" = 4.0
    "F60" = 35.45
    "K60" = 5.0
    "F68" = 27.84
    "K68" = 4.0
    "F89" = 19.24
    "K89" = 2.0
    "F92" = 55.52
    "K92" = 8.0
    "F99" = 69.45
    "K99" = 5.0
    "F110" = 37.45
    "K110" = 1.0
    "F111" = 149.8
    "K111" = 4.0
    "F114" = 224.7
    "K114" = 6.0
    "F121" = 746.86
    "K121" = 2.0
    "F125" = 140.0
    "K125" = 7.0
    "F126" = 53.45
    "K126" = 5.0
    "F130" = 0.0
    "K130" = 0.0
    "F131" = 381.99
    "K131" = 3.0
    "F132" = 11.22
    "K132" = 1.0
    "F133" = 121.8
    "K133" = 12.0
    "F140" = 64.1
    "K140" = 5.0
    "F145" = 37.3
    "K145" = 5.0
    "F147" = 12.72
    "K147" = 4.0
    "F151" = 34.2
    "K151" = 4.0
    "F152" = 57.68
    "K152" = 7.0
    "F164" = 9.17
    "K164" = 7.0
    "F166" = 17.04
    "K166" = 4.0
    "F167" = 150.64
    "K167" = 4.0
    "F168" = 6.96
    "K168" = 3.0
    "F170" = 96.25
    "K170" = 5.0
    "F171" = 10.24
    "K171" = 8.0
    "F173" = 9.8
    "K173" = 4.0
    "F177" = 10.68
    "K177" = 2.0
    "F178" = 58.8
    "K178" = 5.0
    "F183" = 53498.05
    "K183" = 975.0
    "F185" = 53498.05
    "K185" = 975.0
}

foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Value = $cellValues[$ref]
}
